$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2615716457366943

$ws.Range("B3").Value = 0.2606146335601807
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 104.8999229583518
